# Petty cash book - 7-Jun-2021 end of day update
# All transactions for the period 31-May-2021 .. 7-Jun-2021 are cleared out
# (the petty cash book is "rolled forward"): the opening balance becomes
# the new carried-forward balance, the first transaction row is re-dated
# to 7-Jun-2021 and the rest of the rows in that block are emptied.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New carried-forward opening balance
$ws.Range("E2").Value = 942625

# Row 3: re-date to 7-Jun-2021, drop the Debit amount/formula, keep the
# "Wages Expense" label and the running-balance formula in E3.
$ws.Range("A3").Value = 44354
$ws.Range("D3").Clear()

# Rows 4-34: clear every cell except the running-balance formula in column E.
$rowsToClear = @(
    @{ Row = 4;  Cols = @("B","D") },
    @{ Row = 5;  Cols = @("B","C") },
    @{ Row = 6;  Cols = @("B","D") },
    @{ Row = 7;  Cols = @("B","D") },
    @{ Row = 8;  Cols = @("B","C") },
    @{ Row = 9;  Cols = @("B","D") },
    @{ Row = 10; Cols = @("B","D") },
    @{ Row = 11; Cols = @("B","C") },
    @{ Row = 12; Cols = @("A","B","D") },
    @{ Row = 13; Cols = @("B","D") },
    @{ Row = 14; Cols = @("B","D") },
    @{ Row = 15; Cols = @("B","C") },
    @{ Row = 16; Cols = @("B","C") },
    @{ Row = 17; Cols = @("B","D") },
    @{ Row = 18; Cols = @("B","D") },
    @{ Row = 19; Cols = @("A","B","D") },
    @{ Row = 20; Cols = @("B","D") },
    @{ Row = 21; Cols = @("B","C") },
    @{ Row = 22; Cols = @("B","D") },
    @{ Row = 23; Cols = @("B","D") },
    @{ Row = 24; Cols = @("B","C") },
    @{ Row = 25; Cols = @("B","D") },
    @{ Row = 26; Cols = @("A","B","D") },
    @{ Row = 27; Cols = @("B","D") },
    @{ Row = 28; Cols = @("B","D") },
    @{ Row = 29; Cols = @("B","D") },
    @{ Row = 30; Cols = @("B","C") },
    @{ Row = 31; Cols = @("B","C") },
    @{ Row = 32; Cols = @("B","D") },
    @{ Row = 33; Cols = @("B","C") },
    @{ Row = 34; Cols = @("B","D") },
    @{ Row = 35; Cols = @("A") }
)

foreach ($item in $rowsToClear) {
    $r = $item.Row
    foreach ($col in $item.Cols) {
        $ws.Range("$col$r").Clear()
    }
}

# Scroll the frozen pane back up to the top of the new period and select D4
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("D4").Select()
